# Word COM-interop edit script
# Applies the changes described by the commit "fix critical path":
#  1. Title paragraphs (1 & 2): collapse the spell-check-split runs into
#     single runs and tag everything with lang="en-US" (run + paragraph
#     mark formatting).
#  2. The "critical path" paragraph: split the long run describing the
#     critical path into three runs, inserting the new
#     " check specifications ->" text, and move the _GoBack bookmark so it
#     sits between the new runs instead of after "(floats)".

$d = $word.ActiveDocument

# --- 1. Title paragraphs -------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$titleRange = $d.Range($p1.Range.Start, $p2.Range.End)

$titleXml = '<w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Jan Rehwaldt, 2012-02-17, University of Tartu</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Exercise 2, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Project Management</w:t></w:r></w:p>'

$titleRange.InsertXML($titleXml)

# --- 2. "Critical path" paragraph ----------------------------------------
# Locate the paragraph that contains the critical-path sentence (its index
# shifts if the title paragraphs above changed paragraph count, so search
# for it rather than hard-coding the index).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*critical path*") {
        $targetPara = $candidate
        break
    }
}

$criticalXml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t>critical path</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">is marked in </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">bold red </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#8211; specify overall system -&gt; specify module C -&gt;</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> check specifications -&gt;</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> design module A -&gt; code/test module A -&gt; integration/test system.</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Within this path all possible delays</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (floats)</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> are </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t>0</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>' +
    '</w:p>'

$targetPara.Range.InsertXML($criticalXml)

Write-Output "edit complete"
